$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add column N (html_ready_tweets) ---
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("N1").Value = 'html_ready_tweets'

# --- Propagate row formatting down for new rows 7-11 ---
$ws.Range("A6:M6").Copy() | Out-Null
$ws.Range("A7:M11").PasteSpecial(-4122) | Out-Null
$ws.Range("N1").Copy() | Out-Null
$ws.Range("N2:N11").PasteSpecial(-4122) | Out-Null

$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2:N11").PasteSpecial(-4122) | Out-Null

# --- Cell values ---
# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 'RT @SandraSentinel: 8th commandment from God. "Thou Shall Not Steal"
Our voting process has been stolen.
Our right to have the president w…'
$ws.Range("C2").Value = 1339123251994246912.0
$ws.Range("D2").Value = 140
$ws.Range("E2").Value = 44181.34765046297252411
$ws.Range("F2").Value = 'en'
$ws.Range("G2").Value = 1221764091091440128.0
$ws.Range("H2").Value = 'CHUGGER50357510'
$ws.Range("I2").Value = 95
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = 'RT @SandraSentinel: 8th commandment from God. "Thou Shall Not Steal"Our voting process has been stolen.Our right to have the president w… '
$ws.Range("L2").Value = ' 8th commandment from God. "Thou Shall Not Steal"  Our voting process has been stolen. Our right to have the president w…'
$ws.Range("M2").Value = ' 8th commandment from God. "Thou Shall Not Steal"Our voting process has been stolen.Our right to have the president w… '
$ws.Range("N2").Value = 'RT @SandraSentinel: 8th commandment from God. "Thou Shall Not Steal"Our voting process has been stolen.Our right to have the president w… '

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = '@EdRaposo1 Não existe absolutamente nada confirmado. Ao contrário do que a mídia veicula por aí. QUEM DECIDE É O MI… https://t.co/e148rLNWg2'
$ws.Range("C3").Value = 1339123246055230976.0
$ws.Range("D3").Value = 140
$ws.Range("E3").Value = 44181.3476388888884685
$ws.Range("F3").Value = 'pt'
$ws.Range("G3").Value = 966303579081453568.0
$ws.Range("H3").Value = 'TAMOS_AI_38'
$ws.Range("I3").Value = 2214
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = '@ EdRaposo1 There is absolutely nothing confirmed. Contrary to what the media reports there. WHO DECIDES IS MI… https://t.co/e148rLNWg2 '
$ws.Range("L3").Value = ' Não existe absolutamente nada confirmado. Ao contrário do que a mídia veicula por aí. QUEM DECIDE É O MI… '
$ws.Range("M3").Value = ' EdRaposo1 There is absolutely nothing confirmed. Contrary to what the media reports there. WHO DECIDES IS MI…  '
$ws.Range("N3").Value = '@ EdRaposo1 There is absolutely nothing confirmed. Contrary to what the media reports there. WHO DECIDES IS MI…  '

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 'RT @VoicePoliticsmg: A bust water pipe could delay Georgia vote count, no ballots ruined.
#USAElections2020 #Biden #Trump #PresidentialEle…'
$ws.Range("C4").Value = 1339123245291885056.0
$ws.Range("D4").Value = 140
$ws.Range("E4").Value = 44181.3476388888884685
$ws.Range("F4").Value = 'en'
$ws.Range("G4").Value = 217271936
$ws.Range("H4").Value = 'MelissaAtwoodTx'
$ws.Range("I4").Value = 310
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = 'RT @VoicePoliticsmg: A burst water pipe could delay Georgia vote count, no ballots ruined.#USAElections2020 #Biden #Trump #PresidentialEle… '
$ws.Range("L4").Value = ' A bust water pipe could delay Georgia vote count, no ballots ruined.  #USAElections2020 #Biden #Trump #PresidentialEle…'
$ws.Range("M4").Value = ' A burst water pipe could delay Georgia vote count, no ballots ruined.#USAElections2020 #Biden #Trump #PresidentialEle… '
$ws.Range("N4").Value = 'RT @VoicePoliticsmg: A burst water pipe could delay Georgia vote count, no ballots ruined.#USAElections2020 #Biden #Trump #PresidentialEle… '

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'RT @DLFNachrichten: In der Causa #Woelki fordert der Münsteraner Kirchenrechtler Schueller Konsequenzen. Woelki leide offenbar an Realitäts…'
$ws.Range("C5").Value = 1339123205227811072.0
$ws.Range("D5").Value = 140
$ws.Range("E5").Value = 44181.34752314814977581
$ws.Range("F5").Value = 'de'
$ws.Range("G5").Value = 3769471457
$ws.Range("H5").Value = 'smartie1944'
$ws.Range("I5").Value = 44
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = 'RT @DLFNachrichten: In the case of #Woelki, the Muenster canon lawyer Schueller demands consequences. Woelki apparently suffers from reality ... '
$ws.Range("L5").Value = ' In der Causa #Woelki fordert der Münsteraner Kirchenrechtler Schueller Konsequenzen. Woelki leide offenbar an Realitäts…'
$ws.Range("M5").Value = ' In the case of #Woelki, the Muenster canon lawyer Schueller demands consequences. Woelki apparently suffers from reality ... '
$ws.Range("N5").Value = 'RT @DLFNachrichten: In the case of #Woelki, the Muenster canon lawyer Schueller demands consequences. Woelki apparently suffers from reality ... '

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'RT @Trump_Fact_News: Melania #Trump va à l''hôpital pour lire des histoires a des enfants malades, aucun média ne couvre l''événement à part…'
$ws.Range("C6").Value = 1339123197220969984.0
$ws.Range("D6").Value = 139
$ws.Range("E6").Value = 44181.34750000000349246
$ws.Range("F6").Value = 'fr'
$ws.Range("G6").Value = 376695755
$ws.Range("H6").Value = 'WhyAlwaysMeHaha'
$ws.Range("I6").Value = 182
$ws.Range("J6").Value = 'Paris'
$ws.Range("K6").Value = 'RT @Trump_Fact_News: Melania #Trump goes to the hospital to read stories to sick children, no media is covering the event except ... '
$ws.Range("L6").Value = ' Melania #Trump va à l''hôpital pour lire des histoires a des enfants malades, aucun média ne couvre l''événement à part…'
$ws.Range("M6").Value = ' Melania #Trump goes to the hospital to read stories to sick children, no media is covering the event except ... '
$ws.Range("N6").Value = 'RT @Trump_Fact_News: Melania #Trump goes to the hospital to read stories to sick children, no media is covering the event except ... '

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = '📢 Včera se pravděpodobně #Brexit posunul na #Deal stranu. Graf mluví jasně. Do toho se rýsuje #US #Stimulus (po fin… https://t.co/USIcrbI8uO'
$ws.Range("C7").Value = 1339123162261413888.0
$ws.Range("D7").Value = 140
$ws.Range("E7").Value = 44181.3474074073965312
$ws.Range("F7").Value = 'cs'
$ws.Range("G7").Value = 212985092
$ws.Range("H7").Value = 'alan_kooper'
$ws.Range("I7").Value = 429
$ws.Range("J7").Value = 'Czech Republic'
$ws.Range("K7").Value = '📢 Yesterday, #Brexit probably moved to the #Deal side. The graph speaks clearly. #US #Stimulus is outlined (after fin… https://t.co/USIcrbI8uO '
$ws.Range("L7").Value = '📢 Včera se pravděpodobně #Brexit posunul na #Deal stranu. Graf mluví jasně. Do toho se rýsuje #US #Stimulus (po fin… '
$ws.Range("M7").Value = '📢 Yesterday, #Brexit probably moved to the #Deal side. The graph speaks clearly. #US #Stimulus is outlined (after fin…  '
$ws.Range("N7").Value = '📢 Yesterday, #Brexit probably moved to the #Deal side. The graph speaks clearly. #US #Stimulus is outlined (after fin…  '

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'RT @Trump_Fact_News: Donc on voudrait nous faire croire que ce vieil homme sénile que les médias appellent "Président des Etats-unis" a fai…'
$ws.Range("C8").Value = 1339123160621477888.0
$ws.Range("D8").Value = 140
$ws.Range("E8").Value = 44181.3473958333270275
$ws.Range("F8").Value = 'fr'
$ws.Range("G8").Value = 1294368775165424128.0
$ws.Range("H8").Value = 'NoUseFo90037094'
$ws.Range("I8").Value = 81
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = 'RT @Trump_Fact_News: So we would like to make us believe that this senile old man whom the media calls "President of the United States" did… '
$ws.Range("L8").Value = ' Donc on voudrait nous faire croire que ce vieil homme sénile que les médias appellent "Président des Etats-unis" a fai…'
$ws.Range("M8").Value = ' So we would like to make us believe that this senile old man whom the media calls "President of the United States" did… '
$ws.Range("N8").Value = 'RT @Trump_Fact_News: So we would like to make us believe that this senile old man whom the media calls "President of the United States" did… '

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'Microsoft Fixes Xbox Series X 4K Blu-Ray Brightness Bug SEE MORE HERE ==&gt; https://t.co/uTKq3jCQuw #nintendoswitch… https://t.co/LcCuPRHQhP'
$ws.Range("C9").Value = 1339123147484880896.0
$ws.Range("D9").Value = 141
$ws.Range("E9").Value = 44181.34736111111124046
$ws.Range("F9").Value = 'en'
$ws.Range("G9").Value = 306127388
$ws.Range("H9").Value = 'bitcoinconnect'
$ws.Range("I9").Value = 1253
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = 'Microsoft Fixes Xbox Series X 4K Blu-Ray Brightness Bug SEE MORE HERE ==&gt; https://t.co/uTKq3jCQuw #nintendoswitch… https://t.co/LcCuPRHQhP '
$ws.Range("L9").Value = 'Microsoft Fixes Xbox Series X 4K Blu-Ray Brightness Bug SEE MORE HERE ==&gt;  #nintendoswitch… '
$ws.Range("M9").Value = 'Microsoft Fixes Xbox Series X 4K Blu-Ray Brightness Bug SEE MORE HERE ==&gt;  #nintendoswitch…  '
$ws.Range("N9").Value = 'Microsoft Fixes Xbox Series X 4K Blu-Ray Brightness Bug SEE MORE HERE ==&gt;  #nintendoswitch…  '

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'RT @bitcoinconnect: I thought the PS5 DualSense controller was a gimmick — until I played this game SEE MORE HERE ==&gt; https://t.co/8PTOFMzQ…'
$ws.Range("C10").Value = 1339123096993852928.0
$ws.Range("D10").Value = 143
$ws.Range("E10").Value = 44181.34722222221898846
$ws.Range("F10").Value = 'en'
$ws.Range("G10").Value = 1332740896199237888.0
$ws.Range("H10").Value = 'XboxRetweeter'
$ws.Range("I10").Value = 396
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = 'RT @bitcoinconnect: I thought the PS5 DualSense controller was a gimmick — until I played this game SEE MORE HERE ==&gt; https://t.co/8PTOFMzQ… '
$ws.Range("L10").Value = ' I thought the PS5 DualSense controller was a gimmick — until I played this game SEE MORE HERE ==&gt; …'
$ws.Range("M10").Value = ' I thought the PS5 DualSense controller was a gimmick — until I played this game SEE MORE HERE ==&gt; … '
$ws.Range("N10").Value = 'RT @bitcoinconnect: I thought the PS5 DualSense controller was a gimmick — until I played this game SEE MORE HERE ==&gt; … '

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'Don''t be part of the surrender caucus, fight for #America, @senatemajldr. #Trump was cheated, reclaim republican vi… https://t.co/JAJjxBhIFE'
$ws.Range("C11").Value = 1339123056816442880.0
$ws.Range("D11").Value = 140
$ws.Range("E11").Value = 44181.34711805555707542
$ws.Range("F11").Value = 'en'
$ws.Range("G11").Value = 35138538
$ws.Range("H11").Value = 'marcvincens'
$ws.Range("I11").Value = 116
$ws.Range("J11").Value = 'USA'
$ws.Range("K11").Value = 'Don''t be part of the surrender caucus, fight for #America, @senatemajldr. #Trump was cheated, reclaim republican vi… https://t.co/JAJjxBhIFE '
$ws.Range("L11").Value = 'Don''t be part of the surrender caucus, fight for #America, . #Trump was cheated, reclaim republican vi… '
$ws.Range("M11").Value = 'Don''t be part of the surrender caucus, fight for #America, . #Trump was cheated, reclaim republican vi…  '
$ws.Range("N11").Value = 'Don''t be part of the surrender caucus, fight for #America, @senatemajldr. #Trump was cheated, reclaim republican vi…  '

